$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 578.2857
$ws.Range("I9").Value = 507.3846
$ws.Range("J9").Value = 1500
$ws.Range("K9").Value = 507.3846
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = -338.3846
$ws.Range("N9").Value = -1838
$ws.Range("H33").Value = 916.2632
$ws.Range("I33").Value = 1304.6154
$ws.Range("K33").Value = 1304.6154
$ws.Range("M33").Value = -1075.6154
$ws.Range("H58").Value = 267427.06
$ws.Range("I58").Value = 2284.4119
$ws.Range("K58").Value = 6853.2357
$ws.Range("M58").Value = -6703.2357
$ws.Range("H113").Value = 8166.625
$ws.Range("I113").Value = 8166.625
$ws.Range("K113").Value = 8166.625
$ws.Range("M113").Value = -4912.625
$ws.Range("H125").Value = 5933.077
$ws.Range("I125").Value = 5578.875
$ws.Range("K125").Value = 50209.875
$ws.Range("M125").Value = -47749.875
$ws.Range("H132").Value = 3185.5178
$ws.Range("I132").Value = 3225.2546
$ws.Range("K132").Value = 9675.763800000001
$ws.Range("M132").Value = -7145.763800000001
$ws.Range("H135").Value = 3014.9546
$ws.Range("I135").Value = 3074
$ws.Range("J135").Value = 2749.25
$ws.Range("K135").Value = 27666
$ws.Range("L135").Value = 24743.25
$ws.Range("M135").Value = -25131
$ws.Range("N135").Value = -29813.25
$ws.Range("H137").Value = 8098.4033
$ws.Range("I137").Value = 1128.3549
$ws.Range("K137").Value = 3385.0647
$ws.Range("M137").Value = -835.0646999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3028.2927
$ws.Range("I32").Value = 3184.6316
$ws.Range("J32").Value = 1048
$ws.Range("K32").Value = 3184.6316
$ws.Range("L32").Value = 1048
$ws.Range("M32").Value = -2897.6316
$ws.Range("N32").Value = -1622
$ws.Range("H44").Value = 29809.8
$ws.Range("J44").Value = 29809.8
$ws.Range("L44").Value = 29809.8
$ws.Range("N44").Value = -30785.8
$ws.Range("H45").Value = 50304.715
$ws.Range("I45").Value = 64953.688
$ws.Range("K45").Value = 64953.688
$ws.Range("M45").Value = -64576.688
$ws.Range("H55").Value = 30021.143
$ws.Range("I55").Value = 18047.5
$ws.Range("K55").Value = 18047.5
$ws.Range("M55").Value = -17732.5
$ws.Range("H61").Value = 2852.6177
$ws.Range("I61").Value = 2018.56
$ws.Range("J61").Value = 5169.4443
$ws.Range("K61").Value = 2018.56
$ws.Range("L61").Value = 5169.4443
$ws.Range("M61").Value = -1806.56
$ws.Range("N61").Value = -5593.4443
$ws.Range("H88").Value = 5013.4287
$ws.Range("I88").Value = 6051.5
$ws.Range("J88").Value = 3629.3333
$ws.Range("K88").Value = 6051.5
$ws.Range("L88").Value = 3629.3333
$ws.Range("M88").Value = -5645.5
$ws.Range("N88").Value = -4441.3333
$ws.Range("H91").Value = 5013.4287
$ws.Range("I91").Value = 6051.5
$ws.Range("J91").Value = 3629.3333
$ws.Range("K91").Value = 6051.5
$ws.Range("L91").Value = 3629.3333
$ws.Range("M91").Value = -4647.5
$ws.Range("N91").Value = -6437.3333
$ws.Range("H132").Value = 406411.84
$ws.Range("I132").Value = 635565.9
$ws.Range("J132").Value = 16850
$ws.Range("K132").Value = 1906697.7
$ws.Range("L132").Value = 50550
$ws.Range("M132").Value = -1904167.7
$ws.Range("N132").Value = -55610
$ws.Range("H136").Value = 2852.6177
$ws.Range("I136").Value = 2018.56
$ws.Range("J136").Value = 5169.4443
$ws.Range("K136").Value = 6055.68
$ws.Range("L136").Value = 15508.3329
$ws.Range("M136").Value = -3505.68
$ws.Range("N136").Value = -20608.3329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1250
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 12822390
$ws.Range("I105").Value = 934.5
$ws.Range("K105").Value = 934.5
$ws.Range("M105").Value = 812.5
$ws.Range("H132").Value = 49998.8
$ws.Range("J132").Value = 49998.8
$ws.Range("L132").Value = 49998.8
$ws.Range("N132").Value = -60118.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 323043.53
$ws.Range("I31").Value = 771370.4399999999
$ws.Range("J31").Value = 14818.75
$ws.Range("K31").Value = 771370.4399999999
$ws.Range("L31").Value = 14818.75
$ws.Range("M31").Value = -771075.4399999999
$ws.Range("N31").Value = -15408.75
$ws.Range("H34").Value = 323043.53
$ws.Range("I34").Value = 771370.4399999999
$ws.Range("J34").Value = 14818.75
$ws.Range("K34").Value = 771370.4399999999
$ws.Range("L34").Value = 14818.75
$ws.Range("M34").Value = -771168.4399999999
$ws.Range("N34").Value = -15222.75
$ws.Range("H94").Value = 887.7778
$ws.Range("I94").Value = 899.5
$ws.Range("J94").Value = 884.4286
$ws.Range("K94").Value = 899.5
$ws.Range("L94").Value = 884.4286
$ws.Range("M94").Value = -448.5
$ws.Range("N94").Value = -1786.4286
$ws.Range("H99").Value = 6487.6665
$ws.Range("I99").Value = 2695
$ws.Range("J99").Value = 7571.2856
$ws.Range("K99").Value = 2695
$ws.Range("L99").Value = 7571.2856
$ws.Range("M99").Value = -1197
$ws.Range("N99").Value = -10567.2856
$ws.Range("H103").Value = 37499.5
$ws.Range("I103").Value = 37499.5
$ws.Range("K103").Value = 37499.5
$ws.Range("M103").Value = -36327.5
$ws.Range("H126").Value = 6487.6665
$ws.Range("I126").Value = 2695
$ws.Range("J126").Value = 7571.2856
$ws.Range("K126").Value = 8085
$ws.Range("L126").Value = 22713.8568
$ws.Range("M126").Value = -5615
$ws.Range("N126").Value = -27653.8568
$ws.Range("H132").Value = 1399.0333
$ws.Range("I132").Value = 1437.6538
$ws.Range("K132").Value = 4312.9614
$ws.Range("M132").Value = -1782.9614
$ws.Range("H141").Value = 79052
$ws.Range("J141").Value = 79052
$ws.Range("L141").Value = 79052
$ws.Range("N141").Value = -89412

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 50.88889
$ws.Range("I38").Value = 50.88889
$ws.Range("K38").Value = 152.66667
$ws.Range("M38").Value = 194.33333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 28654
$ws.Range("J32").Value = 28654
$ws.Range("L32").Value = 28654
$ws.Range("N32").Value = -29246
$ws.Range("H51").Value = 59663
$ws.Range("J51").Value = 59663
$ws.Range("L51").Value = 59663
$ws.Range("N51").Value = -60681
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H122").Value = 348694
$ws.Range("I122").Value = 427581.53
$ws.Range("K122").Value = 1282744.59
$ws.Range("M122").Value = -1280294.59
$ws.Range("H126").Value = 52641104
$ws.Range("I126").Value = 142860850
$ws.Range("J126").Value = 12917.75
$ws.Range("K126").Value = 428582550
$ws.Range("L126").Value = 38753.25
$ws.Range("M126").Value = -428580080
$ws.Range("N126").Value = -43693.25
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1839.9
$ws.Range("I22").Value = 2437.75
$ws.Range("J22").Value = 1441.3334
$ws.Range("K22").Value = 2437.75
$ws.Range("L22").Value = 1441.3334
$ws.Range("M22").Value = -2142.75
$ws.Range("N22").Value = -2031.3334
$ws.Range("H27").Value = 1839.9
$ws.Range("I27").Value = 2437.75
$ws.Range("J27").Value = 1441.3334
$ws.Range("K27").Value = 2437.75
$ws.Range("L27").Value = 1441.3334
$ws.Range("M27").Value = -2330.75
$ws.Range("N27").Value = -1655.3334
$ws.Range("H40").Value = 1007587.4
$ws.Range("I40").Value = 1435424.1
$ws.Range("K40").Value = 1435424.1
$ws.Range("M40").Value = -1435288.1
$ws.Range("H46").Value = 3381.1
$ws.Range("J46").Value = 4036.4546
$ws.Range("L46").Value = 4036.4546
$ws.Range("N46").Value = -4412.4546
$ws.Range("H55").Value = 66667850
$ws.Range("I55").Value = 724
$ws.Range("K55").Value = 724
$ws.Range("M55").Value = -551

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2693.7
$ws.Range("I81").Value = 2693.7
$ws.Range("K81").Value = 5387.4
$ws.Range("M81").Value = -4326.4
$ws.Range("H84").Value = 2693.7
$ws.Range("I84").Value = 2693.7
$ws.Range("K84").Value = 26937
$ws.Range("M84").Value = -21633
$ws.Range("H106").Value = 49000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 7484.8184
$ws.Range("I126").Value = 2242.0908
$ws.Range("J126").Value = 12727.546
$ws.Range("K126").Value = 6726.2724
$ws.Range("L126").Value = 38182.638
$ws.Range("M126").Value = -4256.2724
$ws.Range("N126").Value = -43122.638
$ws.Range("H132").Value = 25256.918
$ws.Range("I132").Value = 1789.2646
$ws.Range("J132").Value = 78450.266
$ws.Range("K132").Value = 5367.793799999999
$ws.Range("L132").Value = 235350.798
$ws.Range("M132").Value = -2837.793799999999
$ws.Range("N132").Value = -240410.798
$ws.Range("H136").Value = 370913
$ws.Range("I136").Value = 438457.66
$ws.Range("J136").Value = 215560.3
$ws.Range("K136").Value = 1315372.98
$ws.Range("L136").Value = 646680.8999999999
$ws.Range("M136").Value = -1312822.98
$ws.Range("N136").Value = -651780.8999999999

